$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet currently has a 2-row header (row1: units over E/G/I/J/K, row2: Hiver/Ete/Annee
# sub-headers) followed by 13 data rows (rows 3-15). The target layout is a single header
# row (row1: idx, idx2, Name, Date Start, Date End, (m3/s), (MW1), (MW2), (GWh) Winter,
# (GWh) Summer, (GWh) Year) followed by the same 13 data rows shifted up to rows 2-14.

# Remove the old second header row (Hiver/Ete/Annee) - this shifts the 13 data rows up
# from rows 3-15 to rows 2-14.
$ws.Rows.Item(2).Delete()

# Clear the old first header row's contents/formatting so we can rebuild it.
$ws.Rows.Item(1).Clear()

# Write the new single header row.
$ws.Range("A1").Value2 = "idx"
$ws.Range("B1").Value2 = "idx2"
$ws.Range("C1").Value2 = "Name"
$ws.Range("D1").Value2 = "Date Start"
$ws.Range("E1").Value2 = "Date End"
$ws.Range("F1").Value2 = "(m3/s)"
$ws.Range("G1").Value2 = "(MW1)"
$ws.Range("H1").Value2 = "(MW2)"
$ws.Range("I1").Value2 = "(GWh) Winter"
$ws.Range("J1").Value2 = "(GWh) Summer"
$ws.Range("K1").Value2 = "(GWh) Year"

# Give the numeric-column headers (F1:K1) a distinct style: same font as the rest of the
# data table (9pt) but general (non-numeric) formatting, since they hold text.
$hdr = $ws.Range("F1:K1")
$hdr.Font.Size = 9
$hdr.NumberFormat = "General"

# Restore the selection to the first data row, matching the edited file.
$ws.Range("A2:K2").Select()
